$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 66, shifting existing rows 66..161 down to 67..162
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with the new data record
$ws.Range("A66").Value2 = 3
$ws.Range("B66").Value2 = "Femacal de La Calera"
$ws.Range("C66").Value2 = "Coquimbo"
$ws.Range("D66").Value2 = 44757
$ws.Range("D66").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E66").Value2 = 5
$ws.Range("F66").Value2 = 100112026
$ws.Range("G66").Value2 = "Haba"
$ws.Range("H66").Value2 = "Sin especificar"
$ws.Range("I66").Value2 = "Primera"
$ws.Range("J66").Value2 = 40
$ws.Range("K66").Value2 = 21000
$ws.Range("L66").Value2 = 21000
$ws.Range("M66").Value2 = 21000
$ws.Range("N66").Value2 = "`$/saco 25 kilos"
$ws.Range("O66").Value2 = "Provincia de Limarí"
$ws.Range("P66").Value2 = 840
$ws.Range("Q66").Value2 = 25
$ws.Range("R66").Value2 = "Hortaliza"
